# Auto-generated edit script: updates Universalis market-price derived
# columns (H-N) across the leve-profit worksheets, reflecting refreshed
# price data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1000.8947
$ws.Range("I2").Value = 335.66666
$ws.Range("K2").Value = 335.66666
$ws.Range("M2").Value = -222.66666
$ws.Range("H5").Value = 44.75
$ws.Range("I5").Value = 44.75
$ws.Range("K5").Value = 44.75
$ws.Range("M5").Value = 70.25
$ws.Range("H9").Value = 111.4
$ws.Range("I9").Value = 110.57143
$ws.Range("J9").Value = 113.333336
$ws.Range("K9").Value = 110.57143
$ws.Range("L9").Value = 113.333336
$ws.Range("M9").Value = 58.42856999999999
$ws.Range("N9").Value = -451.333336
$ws.Range("H17").Value = 1940
$ws.Range("I17").Value = 1940
$ws.Range("K17").Value = 5820
$ws.Range("M17").Value = -5652
$ws.Range("H52").Value = 1050
$ws.Range("J52").Value = 1050
$ws.Range("L52").Value = 3150
$ws.Range("N52").Value = -3470
$ws.Range("H62").Value = 4483.45
$ws.Range("I62").Value = 5017.933
$ws.Range("K62").Value = 5017.933
$ws.Range("M62").Value = -4393.933
$ws.Range("H65").Value = 4483.45
$ws.Range("I65").Value = 5017.933
$ws.Range("K65").Value = 25089.665
$ws.Range("M65").Value = -21969.665
$ws.Range("H80").Value = 515.4167
$ws.Range("J80").Value = 968.4
$ws.Range("L80").Value = 2905.2
$ws.Range("N80").Value = -4901.2
$ws.Range("H83").Value = 515.4167
$ws.Range("J83").Value = 968.4
$ws.Range("L83").Value = 8715.6
$ws.Range("N83").Value = -18699.6
$ws.Range("H132").Value = 2738.1277
$ws.Range("I132").Value = 2678.05
$ws.Range("K132").Value = 8034.150000000001
$ws.Range("M132").Value = -5504.150000000001
$ws.Range("H138").Value = 2644.4092
$ws.Range("I138").Value = 1256.875
$ws.Range("J138").Value = 3437.2856
$ws.Range("K138").Value = 3770.625
$ws.Range("L138").Value = 10311.8568
$ws.Range("M138").Value = 1369.375
$ws.Range("N138").Value = -20591.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 279.7
$ws.Range("I5").Value = 263.375
$ws.Range("J5").Value = 345
$ws.Range("K5").Value = 263.375
$ws.Range("L5").Value = 345
$ws.Range("M5").Value = -151.375
$ws.Range("N5").Value = -569
$ws.Range("H32").Value = 8536.617
$ws.Range("I32").Value = 7280.1514
$ws.Range("K32").Value = 7280.1514
$ws.Range("M32").Value = -6993.1514
$ws.Range("H33").Value = 4500
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 6000
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = -2671
$ws.Range("N33").Value = -6658
$ws.Range("H37").Value = 24166.666
$ws.Range("H61").Value = 1548.909
$ws.Range("I61").Value = 1602.4
$ws.Range("K61").Value = 1602.4
$ws.Range("M61").Value = -1390.4
$ws.Range("H132").Value = 3540.75
$ws.Range("I132").Value = 3148.158
$ws.Range("K132").Value = 9444.474
$ws.Range("M132").Value = -6914.474
$ws.Range("H136").Value = 1548.909
$ws.Range("I136").Value = 1602.4
$ws.Range("K136").Value = 4807.200000000001
$ws.Range("M136").Value = -2257.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 279.7
$ws.Range("I4").Value = 263.375
$ws.Range("J4").Value = 345
$ws.Range("K4").Value = 263.375
$ws.Range("L4").Value = 345
$ws.Range("M4").Value = -148.375
$ws.Range("N4").Value = -575
$ws.Range("H7").Value = 425
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 425
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 425
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -651
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H86").Value = 4058.1667
$ws.Range("I86").Value = 3969.8
$ws.Range("K86").Value = 3969.8
$ws.Range("M86").Value = -2846.8
$ws.Range("H89").Value = 4058.1667
$ws.Range("I89").Value = 3969.8
$ws.Range("K89").Value = 19849
$ws.Range("M89").Value = -14233
$ws.Range("H134").Value = 11139.6
$ws.Range("I134").Value = 12324.5
$ws.Range("J134").Value = 6400
$ws.Range("K134").Value = 36973.5
$ws.Range("L134").Value = 19200
$ws.Range("M134").Value = -34438.5
$ws.Range("N134").Value = -24270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4463.143
$ws.Range("I31").Value = 3369.25
$ws.Range("J31").Value = 5921.6665
$ws.Range("K31").Value = 3369.25
$ws.Range("L31").Value = 5921.6665
$ws.Range("M31").Value = -3074.25
$ws.Range("N31").Value = -6511.6665
$ws.Range("H34").Value = 4463.143
$ws.Range("I34").Value = 3369.25
$ws.Range("J34").Value = 5921.6665
$ws.Range("K34").Value = 3369.25
$ws.Range("L34").Value = 5921.6665
$ws.Range("M34").Value = -3167.25
$ws.Range("N34").Value = -6325.6665
$ws.Range("H36").Value = 20000
$ws.Range("I36").Value = 20000
$ws.Range("K36").Value = 20000
$ws.Range("M36").Value = -19612
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 20000
$ws.Range("K40").Value = 20000
$ws.Range("M40").Value = -19840
$ws.Range("H58").Value = 3678.4285
$ws.Range("I58").Value = 5742.7144
$ws.Range("J58").Value = 1614.1428
$ws.Range("K58").Value = 5742.7144
$ws.Range("L58").Value = 1614.1428
$ws.Range("M58").Value = -5539.7144
$ws.Range("N58").Value = -2020.1428
$ws.Range("H134").Value = 1532.375
$ws.Range("I134").Value = 1661.7368
$ws.Range("J134").Value = 1040.8
$ws.Range("K134").Value = 4985.2104
$ws.Range("L134").Value = 3122.4
$ws.Range("M134").Value = -2450.2104
$ws.Range("N134").Value = -8192.4
$ws.Range("H136").Value = 3678.4285
$ws.Range("I136").Value = 5742.7144
$ws.Range("J136").Value = 1614.1428
$ws.Range("K136").Value = 17228.1432
$ws.Range("L136").Value = 4842.428400000001
$ws.Range("M136").Value = -14678.1432
$ws.Range("N136").Value = -9942.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 541.06665
$ws.Range("I5").Value = 508.5
$ws.Range("J5").Value = 997
$ws.Range("K5").Value = 1525.5
$ws.Range("L5").Value = 2991
$ws.Range("M5").Value = -1413.5
$ws.Range("N5").Value = -3215
$ws.Range("H10").Value = 19.9
$ws.Range("I10").Value = 19.9
$ws.Range("K10").Value = 59.7
$ws.Range("M10").Value = 79.30000000000001
$ws.Range("H36").Value = 200
$ws.Range("I36").Value = 200
$ws.Range("K36").Value = 600
$ws.Range("M36").Value = -431
$ws.Range("H104").Value = 29666.666
$ws.Range("I104").Value = 20000
$ws.Range("K104").Value = 60000
$ws.Range("M104").Value = -57379
$ws.Range("H135").Value = 541.06665
$ws.Range("I135").Value = 508.5
$ws.Range("J135").Value = 997
$ws.Range("K135").Value = 4576.5
$ws.Range("L135").Value = 8973
$ws.Range("M135").Value = -2041.5
$ws.Range("N135").Value = -14043

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 857.5
$ws.Range("J13").Value = 857.5
$ws.Range("L13").Value = 857.5
$ws.Range("N13").Value = -1135.5
$ws.Range("H132").Value = 1375
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 816.0741
$ws.Range("I22").Value = 583.8889
$ws.Range("J22").Value = 1280.4445
$ws.Range("K22").Value = 583.8889
$ws.Range("L22").Value = 1280.4445
$ws.Range("M22").Value = -288.8889
$ws.Range("N22").Value = -1870.4445
$ws.Range("H27").Value = 816.0741
$ws.Range("I27").Value = 583.8889
$ws.Range("J27").Value = 1280.4445
$ws.Range("K27").Value = 583.8889
$ws.Range("L27").Value = 1280.4445
$ws.Range("M27").Value = -476.8889
$ws.Range("N27").Value = -1494.4445
$ws.Range("H46").Value = 1552.1538
$ws.Range("I46").Value = 899
$ws.Range("J46").Value = 3729.3333
$ws.Range("K46").Value = 899
$ws.Range("L46").Value = 3729.3333
$ws.Range("M46").Value = -711
$ws.Range("N46").Value = -4105.3333
$ws.Range("H132").Value = 13039.25
$ws.Range("I132").Value = 15364.737
$ws.Range("K132").Value = 46094.211
$ws.Range("M132").Value = -43564.211
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -75060
$ws.Range("H136").Value = 3580.4
$ws.Range("I136").Value = 3477.6667
$ws.Range("K136").Value = 10433.0001
$ws.Range("M136").Value = -7883.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3711.8462
$ws.Range("I136").Value = 3711.8462
$ws.Range("K136").Value = 11135.5386
$ws.Range("M136").Value = -8585.5386
